# B6-PowerPoint.pptx edit
#
# 1) The three data tables (originally on slides 14-16) get re-styled from
#    the deck's custom "Table_0" style to the built-in
#    "Medium Style 2 - Accent 1" table style.
# 2) The presentation's theme colour scheme is repointed from the
#    "Integral / Red Violet" palette to the standard "Office" palette.

$p = $ppt.ActivePresentation

$oldTableStyle = "{6F142DA8-6B31-4C6D-906B-4FFC6825DAC3}"
$newTableStyle = "{2888D82F-CBC0-4233-BDEC-24F373437417}"

for ($sn = 1; $sn -le $p.Slides.Count; $sn++) {
    $s = $p.Slides.Item($sn)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            if ($tbl.Style -eq $oldTableStyle) {
                $tbl.ApplyStyle($newTableStyle)
            }
        }
    }
}

# Re-colour the theme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) to the
# stock Office palette.
$officeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $rgbHex = $officeColors[$i - 1]
    $r = ($rgbHex -shr 16) -band 0xFF
    $g = ($rgbHex -shr 8) -band 0xFF
    $b = $rgbHex -band 0xFF
    # COM RGB() packs colours as 0x00BBGGRR.
    $tcs.Item($i).RGB = $r + ($g * 256) + ($b * 65536)
}
